# Update the group name placeholder from "96Group" to "6014Group"
# on both the Input and Output sheets, then leave the selection on
# the Output sheet at B8 (mirroring the saved cursor position).

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B11").Value = "6014Group"

$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Range("B1").Value = "6014Group"

$wsOutput.Activate()
$wsOutput.Range("B8").Select()
